$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2000166.2
$ws.Range("I6").Value = 2500187.8
$ws.Range("K6").Value = 7500563.399999999
$ws.Range("M6").Value = -7500451.399999999
$ws.Range("H11").Value = 221.5
$ws.Range("I11").Value = 221.5
$ws.Range("K11").Value = 221.5
$ws.Range("M11").Value = -81.5
$ws.Range("H99").Value = 517.3333
$ws.Range("I99").Value = 517.3333
$ws.Range("K99").Value = 1551.9999
$ws.Range("M99").Value = -53.99990000000003
$ws.Range("H100").Value = 1878.2222
$ws.Range("I100").Value = 1822.9375
$ws.Range("K100").Value = 1822.9375
$ws.Range("M100").Value = -1281.9375
$ws.Range("H116").Value = 6333.3687
$ws.Range("I116").Value = 5535.6
$ws.Range("K116").Value = 5535.6
$ws.Range("M116").Value = -2093.6
$ws.Range("H132").Value = 1696.0588
$ws.Range("I132").Value = 1734.0212
$ws.Range("J132").Value = 1250
$ws.Range("K132").Value = 5202.063599999999
$ws.Range("L132").Value = 3750
$ws.Range("M132").Value = -2672.063599999999
$ws.Range("N132").Value = -8810

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7653.6963
$ws.Range("I32").Value = 3636.6567
$ws.Range("K32").Value = 3636.6567
$ws.Range("M32").Value = -3349.6567
$ws.Range("H63").Value = 3250
$ws.Range("J63").Value = 3250
$ws.Range("L63").Value = 3250
$ws.Range("N63").Value = -4622
$ws.Range("H66").Value = 3250
$ws.Range("J66").Value = 3250
$ws.Range("L66").Value = 16250
$ws.Range("N66").Value = -23114
$ws.Range("H88").Value = 2071.4285
$ws.Range("I88").Value = 1900
$ws.Range("J88").Value = 2500
$ws.Range("K88").Value = 1900
$ws.Range("L88").Value = 2500
$ws.Range("M88").Value = -1494
$ws.Range("N88").Value = -3312
$ws.Range("H91").Value = 2071.4285
$ws.Range("I91").Value = 1900
$ws.Range("J91").Value = 2500
$ws.Range("K91").Value = 1900
$ws.Range("L91").Value = 2500
$ws.Range("M91").Value = -496
$ws.Range("N91").Value = -5308
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H134").Value = 150000
$ws.Range("J134").Value = 150000
$ws.Range("L134").Value = 150000
$ws.Range("N134").Value = -160140

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 77404.19
$ws.Range("I99").Value = 49119.668
$ws.Range("J99").Value = 176400
$ws.Range("K99").Value = 49119.668
$ws.Range("L99").Value = 176400
$ws.Range("M99").Value = -47621.668
$ws.Range("N99").Value = -179396
$ws.Range("H107").Value = 1326.4
$ws.Range("I107").Value = 1188.3529
$ws.Range("J107").Value = 2108.6667
$ws.Range("K107").Value = 1188.3529
$ws.Range("L107").Value = 2108.6667
$ws.Range("M107").Value = 731.6470999999999
$ws.Range("N107").Value = -5948.6667
$ws.Range("H117").Value = 45000
$ws.Range("J117").Value = 45000
$ws.Range("L117").Value = 45000
$ws.Range("N117").Value = -54178
$ws.Range("H119").Value = 44999.5
$ws.Range("J119").Value = 44999.5
$ws.Range("L119").Value = 44999.5
$ws.Range("N119").Value = -54675.5
$ws.Range("H134").Value = 2115.25
$ws.Range("I134").Value = 1421.8948
$ws.Range("J134").Value = 4750
$ws.Range("K134").Value = 4265.6844
$ws.Range("L134").Value = 14250
$ws.Range("M134").Value = -1730.6844
$ws.Range("N134").Value = -19320
$ws.Range("H140").Value = 85390
$ws.Range("J140").Value = 85390
$ws.Range("L140").Value = 85390
$ws.Range("N140").Value = -95750

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H118").Value = 250000
$ws.Range("J118").Value = 250000
$ws.Range("L118").Value = 250000
$ws.Range("N118").Value = -253314
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("M121").ClearContents()
$ws.Range("H122").Value = 1985.6875
$ws.Range("I122").Value = 888.36365
$ws.Range("K122").Value = 2665.09095
$ws.Range("M122").Value = -215.0909499999998
$ws.Range("H134").Value = 11145.296
$ws.Range("I134").Value = 7351
$ws.Range("K134").Value = 22053
$ws.Range("M134").Value = -19518

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 1454.9231
$ws.Range("J41").Value = 2607
$ws.Range("L41").Value = 7821
$ws.Range("N41").Value = -8497
$ws.Range("H51").Value = 999
$ws.Range("I51").Value = 999
$ws.Range("K51").Value = 2997
$ws.Range("M51").Value = -2537
$ws.Range("H87").Value = 12565.818
$ws.Range("I87").Value = 8091.0586
$ws.Range("K87").Value = 24273.1758
$ws.Range("M87").Value = -23025.1758
$ws.Range("H90").Value = 12565.818
$ws.Range("I90").Value = 8091.0586
$ws.Range("K90").Value = 72819.52740000001
$ws.Range("M90").Value = -66579.52740000001
$ws.Range("H141").Value = 53400.684
$ws.Range("I141").Value = 968.8570999999999
$ws.Range("K141").Value = 2906.5713
$ws.Range("M141").Value = 2273.4287

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2250
$ws.Range("I80").Value = 2000
$ws.Range("K80").Value = 2000
$ws.Range("M80").Value = -1002
$ws.Range("H83").Value = 2250
$ws.Range("I83").Value = 2000
$ws.Range("K83").Value = 10000
$ws.Range("M83").Value = -5008
$ws.Range("H102").Value = 50000850
$ws.Range("I102").Value = 877.73334
$ws.Range("K102").Value = 877.73334
$ws.Range("M102").Value = 744.26666
$ws.Range("H136").Value = 24847.576
$ws.Range("J136").Value = 24847.576
$ws.Range("L136").Value = 74542.728
$ws.Range("N136").Value = -79642.728

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3118
$ws.Range("I7").Value = 2967.6956
$ws.Range("K7").Value = 2967.6956
$ws.Range("M7").Value = -2855.6956
$ws.Range("H16").Value = 1748.0698
$ws.Range("I16").Value = 1349.5405
$ws.Range("K16").Value = 1349.5405
$ws.Range("M16").Value = -1179.5405
$ws.Range("H22").Value = 2056.1924
$ws.Range("I22").Value = 1486.1666
$ws.Range("J22").Value = 2544.7856
$ws.Range("K22").Value = 1486.1666
$ws.Range("L22").Value = 2544.7856
$ws.Range("M22").Value = -1191.1666
$ws.Range("N22").Value = -3134.7856
$ws.Range("H27").Value = 2056.1924
$ws.Range("I27").Value = 1486.1666
$ws.Range("J27").Value = 2544.7856
$ws.Range("K27").Value = 1486.1666
$ws.Range("L27").Value = 2544.7856
$ws.Range("M27").Value = -1379.1666
$ws.Range("N27").Value = -2758.7856
$ws.Range("H40").Value = 6168.5835
$ws.Range("I40").Value = 5854
$ws.Range("J40").Value = 6609
$ws.Range("K40").Value = 5854
$ws.Range("L40").Value = 6609
$ws.Range("M40").Value = -5718
$ws.Range("N40").Value = -6881
$ws.Range("H102").Value = 97999.336
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 97999.336
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 97999.336
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -104489.336
$ws.Range("H126").Value = 3118
$ws.Range("I126").Value = 2967.6956
$ws.Range("K126").Value = 8903.086800000001
$ws.Range("M126").Value = -6433.086800000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 19500
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 19500
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H100").Value = 875.44446
$ws.Range("I100").Value = 862.1539
$ws.Range("J100").Value = 910
$ws.Range("K100").Value = 1724.3078
$ws.Range("L100").Value = 1820
$ws.Range("M100").Value = -1183.3078
$ws.Range("N100").Value = -2902
$ws.Range("H122").Value = 2005
$ws.Range("I122").Value = 1919.091
$ws.Range("J122").Value = 2950
$ws.Range("K122").Value = 5757.272999999999
$ws.Range("L122").Value = 8850
$ws.Range("M122").Value = -3307.272999999999
$ws.Range("N122").Value = -13750
$ws.Range("H132").Value = 2109.15
$ws.Range("I132").Value = 1876.3529
$ws.Range("K132").Value = 5629.0587
$ws.Range("M132").Value = -3099.0587
